# Rename the "4.3.2 (Jan 26th)" sheet to "4.3.2 (Jan 22nd)"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("4.3.2 (Jan 26th)")
$ws1.Name = "4.3.2 (Jan 22nd)"

# Make this sheet the active/selected tab, with D27 as the current selection
# (this also clears the previous tab selection on the other sheet, which
# keeps its own selection at E24)
$ws1.Activate() | Out-Null
$ws1.Range("D27").Select() | Out-Null
